$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date column stays text (matches existing column formatting)
$ws.Range("B208:B257").NumberFormat = "@"

$problems = @(
    "almostunionfind",
    "alphabetanimals",
    "chopwood",
    "fenwick",
    "forests",
    "more10",
    "moviecollection",
    "railroad",
    "supercomputer",
    "turbo",
    "unionfind",
    "virtualfriends",
    "worstweather",
    "committeeassignment",
    "dictionaryattack",
    "diversecontest",
    "expeditiouscubing",
    "firefly",
    "foolingaround",
    "jumpingyoshi",
    "pie",
    "sendmoremoney",
    "taxing",
    "paintings",
    "anti11",
    "aspenavenue",
    "bobby",
    "brevoptimering",
    "busticket",
    "checkingforcorrectness",
    "classrooms",
    "drinkresponsibly",
    "inflagrantedelicto",
    "linearrecurrence",
    "minimumscalar",
    "nesteddolls",
    "ntnuorienteering",
    "sellingspatulas",
    "silverstarstandsalone",
    "ticketpricing",
    "watersheds",
    "ballotboxes",
    "canonical",
    "exam",
    "frustratedqueue",
    "money",
    "planetaris",
    "scenes",
    "stack",
    "zapis"
)

$dates = @(
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-01-26",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-02",
    "2026-02-07",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-08",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09",
    "2026-02-09"
)

$startRow = 208
for ($i = 0; $i -lt $problems.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $problems[$i]
    $ws.Cells.Item($r, 2).Value = $dates[$i]
}
